# Add a new work-log entry on row 19 of the "Taul1" sheet:
#   - Date (A19), start time (B19), end time (C19)
#   - Description string (E19) referencing the new shared string
#     "Skill Increase + Character Advancement"
# The Tunnit (D), accumulated-hours (F) and remaining-hours (G) columns are
# driven by existing shared formulas, so they recalculate automatically.
# Finally, move the active selection to E19 to match the saved workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Paivamaara (date) - 2020-09-11
$ws.Range("A19").Value = 44085
$ws.Range("A19").NumberFormat = $ws.Range("A18").NumberFormat

# Aloitus (start time) - 06:00
$ws.Range("B19").Value = 0.25
$ws.Range("B19").NumberFormat = $ws.Range("B18").NumberFormat

# Lopetus (end time) - 08:00
$ws.Range("C19").Value = 0.33333333333333331
$ws.Range("C19").NumberFormat = $ws.Range("C18").NumberFormat

# Tyoskentely (work description)
$ws.Range("E19").Value = "Skill Increase + Character Advancement"

# Match the saved selection in the workbook
$ws.Range("E19").Select()
